$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.087.31"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "1.892.62"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5024"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3896"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09236"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.382"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").Value = "1.896.86"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.300"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06638"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.227"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").Value = "28.153.60"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D26").Value = "2.113.37"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.548"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.076"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.604"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.608"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.497"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06606"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.340"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02411"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2199"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.218"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6457"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.954"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6064"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.694"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.193"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
